$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of location data appended after the existing last row (241).
$newRows = @(
    @("LIVE, SEA, BEACH", "12.085394745972502, -68.89842911649694", "LionsDive Beach Resort | Mambo Beach | Curaçao", "Willemstad", "Curacao", "loHbMM9JfCs"),
    @("LIVE, SEA, BEACH", "11.98591047553957, -68.64469205447985", "Klein Curaçao | Mermaid Boat Trips", "Klein Curacao", "Curacao", "0ImA9IcyQwA"),
    @("LIVE, SEA, BRDIGE", "12.106351598725128, -68.93525653244262", "Handelskade & Brionplein | Curaçao", "Willemstad", "Curacao", "28U-t3fA9ks"),
    @("LIVE, SEA, BEACH, AIRPORT", "18.03923594493922, -63.120373366239846", "World famous Plane Beach ✈️ - Maho Beach SXM", "Simpson Bay", "Sint Maarten", "LtzkkAeW_Qg")
)

$startRow = 242
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
    $ws.Range("F$r").Value = $row[5]

    # Match the look of the preceding rows: columns A (Category) and E
    # (Country) carry the bordered "s=8" style used throughout the table.
    $ws.Range("A241:F241").Copy()
    $ws.Range("A$r`:F$r").PasteSpecial(-4122)
}

$excel.CutCopyMode = $false

# Move the view roughly to where it ended up after the edit (scrolled down,
# with the next empty row selected).
$win = $excel.ActiveWindow
$win.ScrollRow = 224
$ws.Range("A246").Select()
